# Node editor / New Story + Change Profile edits:
#  1. Remove the _GoBack bookmark that currently sits after "Error Page".
#  2. Delete the whole "Ev. Story Edit Button ändern" bullet paragraph.
#  3. Re-place the _GoBack bookmark at the end of the
#     "Ev. mehr Doodels, Scribbles,..." paragraph (right after its run,
#     i.e. collapsed immediately before that paragraph's mark).

$d = $word.ActiveDocument

function Find-ParagraphByText($text) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.TrimEnd("`r") -eq $text) {
            return $p
        }
    }
    return $null
}

# --- 1. Drop the stray _GoBack bookmark wherever it currently lives ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 2. Delete the "Ev. Story Edit Button ändern" paragraph entirely ---
$storyEditPara = Find-ParagraphByText("Ev. Story Edit Button ändern")
if ($storyEditPara -ne $null) {
    $storyEditPara.Range.Delete()
}

# --- 3. Add the _GoBack bookmark back, collapsed right after the text of
#        "Ev. mehr Doodels, Scribbles,..." (before its paragraph mark).
#        A collapsed range sitting exactly on the paragraph-mark boundary
#        is mis-anchored by Bookmarks.Add, so nudge past it by briefly
#        inserting a placeholder character, bookmarking that, then
#        deleting the placeholder again (the bookmark collapses back to
#        the now-correct location and survives the deletion). ---
$doodelsPara = Find-ParagraphByText("Ev. mehr Doodels, Scribbles,...")
if ($doodelsPara -ne $null) {
    $endRng = $doodelsPara.Range
    $endRng.MoveEnd(1, -1) | Out-Null
    $insertPoint = $endRng.Duplicate
    $insertPoint.Collapse(0)
    $insertPoint.InsertAfter("X")

    $placeholderRng = $d.Range($insertPoint.Start, $insertPoint.End)
    $d.Bookmarks.Add("_GoBack", $placeholderRng)

    $placeholderRng2 = $d.Range($insertPoint.Start, $insertPoint.End)
    $placeholderRng2.Delete()
}
